$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 16 updates (Pre_SBLR summary row) ---
$ws.Range("D16").Value = 171.8
$ws.Range("E16").NumberFormat = "0.0"
$ws.Range("E16").Value = 182
$ws.Range("F16").Value = 51.3
$ws.Range("G16").Value = 611.79999999999995
$ws.Range("H16").Value = 384.7
$ws.Range("I16").Value = 16458.7
$ws.Range("J16").Value = 346.6
$ws.Range("K16").Value = 57.3
$ws.Range("L16").Value = 49.7

# --- Row 17 updates ---
$ws.Range("D17").Value = 117.1
$ws.Range("E17").Value = 157.5
$ws.Range("F17").Value = 38.5
$ws.Range("G17").Value = 229.5
$ws.Range("H17").Value = 161.1
$ws.Range("I17").Value = 108.6
$ws.Range("J17").Value = 287.5
$ws.Range("K17").NumberFormat = "0.0"
$ws.Range("K17").Value = 35

# --- Row 18 updates ---
$ws.Range("D18").Value = 286.2
$ws.Range("E18").Value = 212.1
$ws.Range("F18").Value = 69.900000000000006
$ws.Range("G18").Value = "Infinite "
$ws.Range("J18").Value = 427.4
$ws.Range("K18").Value = 112.3
$ws.Range("L18").Value = 402.8

# --- Selection change (Excel records the last active cell on save) ---
$ws.Range("N13").Select()
